$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.576.01'
$ws.Range("E2").Value = '  -0.57%  '
$ws.Range("D3").Value = '1.860.53'
$ws.Range("E3").Value = '  -0.91%  '
$ws.Range("D4").Value = '''1.011'
$ws.Range("E4").Value = '  +0.71%  '
$ws.Range("D5").Value = '''333.57'
$ws.Range("E5").Value = '  +0.28%  '
$ws.Range("E6").Value = '  +0.67%  '
$ws.Range("D7").Value = '''0.4673'
$ws.Range("E7").Value = '  -0.90%  '
$ws.Range("D8").Value = '''0.3888'
$ws.Range("E8").Value = '  -1.71%  '
$ws.Range("D9").Value = '''45.71'
$ws.Range("D10").Value = '''0.07962'
$ws.Range("E10").Value = '  -1.16%  '
$ws.Range("D11").Value = '''1.001'
$ws.Range("E11").Value = '  -3.02%  '
$ws.Range("D12").Value = '''21.60'
$ws.Range("E12").Value = '  -2.71%  '
$ws.Range("D13").Value = '1.858.54'
$ws.Range("E13").Value = '  -1.15%  '
$ws.Range("D14").Value = '''5.977'
$ws.Range("E14").Value = '  -0.04%  '
$ws.Range("D15").Value = '''7.208'
$ws.Range("E15").Value = '  +1.17%  '
$ws.Range("D16").Value = '''1.012'
$ws.Range("E16").Value = '  +0.78%  '
$ws.Range("D17").Value = '''87.76'
$ws.Range("E17").Value = '  +0.78%  '
$ws.Range("D18").Value = '''0.06706'
$ws.Range("E18").Value = '  +0.55%  '
$ws.Range("E19").Value = '  -0.50%  '
$ws.Range("D20").Value = '''16.90'
$ws.Range("E20").Value = '  -1.90%  '
$ws.Range("E21").Value = '  +0.66%  '
$ws.Range("D22").Value = '27.550.18'
$ws.Range("E22").Value = '  -0.74%  '
$ws.Range("D23").Value = '''5.436'
$ws.Range("E23").Value = '  -1.67%  '
$ws.Range("D24").Value = '''10.84'
$ws.Range("E24").Value = '  -1.48%  '
$ws.Range("D25").Value = '''2.304'
$ws.Range("E25").Value = '  -0.09%  '
$ws.Range("D26").Value = '2.081.77'
$ws.Range("E26").Value = '  -1.03%  '
$ws.Range("D27").Value = '''158.66'
$ws.Range("E27").Value = '  -0.54%  '
$ws.Range("D28").Value = '''19.70'
$ws.Range("E28").Value = '  -2.41%  '
$ws.Range("D29").Value = '''2.121'
$ws.Range("E29").Value = '  +0.91%  '
$ws.Range("D30").Value = '''5.375'
$ws.Range("E30").Value = '  -3.61%  '
$ws.Range("D31").Value = '''121.11'
$ws.Range("E31").Value = '  -0.52%  '
$ws.Range("D32").Value = '''0.9705'
$ws.Range("E32").Value = '  -1.31%  '
$ws.Range("D33").Value = '''0.09448'
$ws.Range("E33").Value = '  -0.82%  '
$ws.Range("D34").Value = '''3.644'
$ws.Range("E34").Value = '  +1.26%  '
$ws.Range("D35").Value = '''5.286'
$ws.Range("E35").Value = '  -1.35%  '
$ws.Range("D36").Value = '''1.325'
$ws.Range("E36").Value = '  -8.37%  '
$ws.Range("D37").Value = '''0.06023'
$ws.Range("E37").Value = '  -1.65%  '
$ws.Range("D38").Value = '''0.02211'
$ws.Range("D39").Value = '''1.192'
$ws.Range("E39").Value = '  -3.21%  '
$ws.Range("D40").Value = '''8.192'
$ws.Range("E40").Value = '  +0.55%  '
$ws.Range("E41").Value = '  +0.67%  '
$ws.Range("D42").Value = '''0.5905'
$ws.Range("E42").Value = '  -1.94%  '
$ws.Range("D43").Value = '''0.1874'
$ws.Range("E43").Value = '  -1.56%  '
$ws.Range("D44").Value = '''10.20'
$ws.Range("E44").Value = '  -0.73%  '
$ws.Range("D45").Value = '''1.250'
$ws.Range("E45").Value = '  -1.44%  '
$ws.Range("D46").Value = '''0.5611'
$ws.Range("E46").Value = '  -2.04%  '
$ws.Range("D47").Value = '''12.06'
$ws.Range("E47").Value = '  -1.51%  '
$ws.Range("D48").Value = '''1.910'
$ws.Range("E48").Value = '  -2.06%  '
$ws.Range("D49").Value = '''3.271'
$ws.Range("E49").Value = '  -3.46%  '
$ws.Range("E50").Value = '  -2.41%  '
$ws.Range("D51").Value = '''112.43'
$ws.Range("E51").Value = '  -1.55%  '
